$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Done" status (with the same green fill formatting used for the other
# API test rows) for the newly added Contacts API test cases in D25:D29.
$ws.Range("D24").Copy()
$ws.Range("D25:D29").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D25:D29").Value = "Done"

# Move the active selection to D30, matching where editing left off.
$ws.Range("D30").Select()
